$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Copy()
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = ""

for ($c = 0; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $c
}
